$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sync boss (stage1) attack range values in column F (RotationSpeed/AttackRange column) to 500
$ws.Range("F2").Value = 500
$ws.Range("F3").Value = 500
$ws.Range("F4").Value = 500
$ws.Range("F5").Value = 500
$ws.Range("F6").Value = 500
$ws.Range("F7").Value = 500
$ws.Range("F8").Value = 500
$ws.Range("F9").Value = 500
$ws.Range("F10").Value = 500
$ws.Range("F11").Value = 500
$ws.Range("F12").Value = 500

# Update active cell selection to match authored state
$ws.Range("L17").Select()
